# Scheduled-runner Sheets update: refresh Leve profit calculations (currentAveragePrice-derived
# columns H-N) across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets with latest market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 925
$ws.Range("J70").Value = 925
$ws.Range("L70").Value = 2775
$ws.Range("N70").Value = -3315
$ws.Range("H73").Value = 925
$ws.Range("J73").Value = 925
$ws.Range("L73").Value = 2775
$ws.Range("N73").Value = -4647
$ws.Range("H96").Value = 1959.5
$ws.Range("I96").Value = 1959.5
$ws.Range("K96").Value = 5878.5
$ws.Range("M96").Value = -4505.5
$ws.Range("H99").Value = 706.86664
$ws.Range("J99").Value = 670.25
$ws.Range("L99").Value = 2010.75
$ws.Range("N99").Value = -5006.75
$ws.Range("H132").Value = 3369248
$ws.Range("I132").Value = 5052061.5
$ws.Range("J132").Value = 3621.111
$ws.Range("K132").Value = 15156184.5
$ws.Range("L132").Value = 10863.333
$ws.Range("M132").Value = -15153654.5
$ws.Range("N132").Value = -15923.333
$ws.Range("H137").Value = 1018.0541
$ws.Range("I137").Value = 844.7925
$ws.Range("J137").Value = 1455.3334
$ws.Range("K137").Value = 2534.3775
$ws.Range("L137").Value = 4366.0002
$ws.Range("M137").Value = 15.62249999999995
$ws.Range("N137").Value = -9466.0002
$ws.Range("H138").Value = 1714.804
$ws.Range("I138").Value = 978.95123
$ws.Range("J138").Value = 4731.8
$ws.Range("K138").Value = 2936.85369
$ws.Range("L138").Value = 14195.4
$ws.Range("M138").Value = 2203.14631
$ws.Range("N138").Value = -24475.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 3932.389
$ws.Range("I41").Value = 2198.7856
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 2198.7856
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = -1784.7856
$ws.Range("N41").Value = -10828
$ws.Range("H75").Value = 29998
$ws.Range("J75").Value = 29998
$ws.Range("L75").Value = 29998
$ws.Range("N75").Value = -31746
$ws.Range("H78").Value = 29998
$ws.Range("J78").Value = 29998
$ws.Range("L78").Value = 89994
$ws.Range("N78").Value = -98730
$ws.Range("H97").Value = 474.32257
$ws.Range("I97").Value = 487.0345
$ws.Range("J97").Value = 290
$ws.Range("K97").Value = 487.0345
$ws.Range("L97").Value = 290
$ws.Range("M97").Value = 8.96550000000002
$ws.Range("N97").Value = -1282
$ws.Range("H132").Value = 1513.8723
$ws.Range("I132").Value = 795.19446
$ws.Range("J132").Value = 3865.9092
$ws.Range("K132").Value = 2385.58338
$ws.Range("L132").Value = 11597.7276
$ws.Range("M132").Value = 144.41662
$ws.Range("N132").Value = -16657.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 9533.333000000001
$ws.Range("I82").Value = 4400
$ws.Range("J82").Value = 12100
$ws.Range("K82").Value = 4400
$ws.Range("L82").Value = 12100
$ws.Range("M82").Value = -4017
$ws.Range("N82").Value = -12866
$ws.Range("H85").Value = 9533.333000000001
$ws.Range("I85").Value = 4400
$ws.Range("J85").Value = 12100
$ws.Range("K85").Value = 4400
$ws.Range("L85").Value = 12100
$ws.Range("M85").Value = -3074
$ws.Range("N85").Value = -14752
$ws.Range("H94").Value = 896.85187
$ws.Range("I94").Value = 443.57144
$ws.Range("K94").Value = 443.57144
$ws.Range("M94").Value = 7.428560000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1769.0151
$ws.Range("I31").Value = 1126.6852
$ws.Range("J31").Value = 4659.5
$ws.Range("K31").Value = 1126.6852
$ws.Range("L31").Value = 4659.5
$ws.Range("M31").Value = -831.6851999999999
$ws.Range("N31").Value = -5249.5
$ws.Range("H34").Value = 1769.0151
$ws.Range("I34").Value = 1126.6852
$ws.Range("J34").Value = 4659.5
$ws.Range("K34").Value = 1126.6852
$ws.Range("L34").Value = 4659.5
$ws.Range("M34").Value = -924.6851999999999
$ws.Range("N34").Value = -5063.5
$ws.Range("H105").Value = 6330.143
$ws.Range("I105").Value = 10766.667
$ws.Range("J105").Value = 3002.75
$ws.Range("K105").Value = 10766.667
$ws.Range("L105").Value = 3002.75
$ws.Range("M105").Value = -9019.666999999999
$ws.Range("N105").Value = -6496.75
$ws.Range("H134").Value = 1900.2858
$ws.Range("I134").Value = 1688.3334
$ws.Range("J134").Value = 2281.8
$ws.Range("K134").Value = 5065.0002
$ws.Range("L134").Value = 6845.400000000001
$ws.Range("M134").Value = -2530.0002
$ws.Range("N134").Value = -11915.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 497
$ws.Range("I5").Value = 411.30435
$ws.Range("K5").Value = 1233.91305
$ws.Range("M5").Value = -1121.91305
$ws.Range("H98").Value = 315.91666
$ws.Range("I98").Value = 249.75
$ws.Range("J98").Value = 349
$ws.Range("K98").Value = 749.25
$ws.Range("L98").Value = 1047
$ws.Range("M98").Value = 748.75
$ws.Range("N98").Value = -4043
$ws.Range("H108").Value = 897.13336
$ws.Range("I108").Value = 454.75
$ws.Range("K108").Value = 1364.25
$ws.Range("M108").Value = 1515.75
$ws.Range("H122").Value = 1000573.2
$ws.Range("I122").Value = 618.05884
$ws.Range("J122").Value = 1515701.5
$ws.Range("K122").Value = 5562.52956
$ws.Range("L122").Value = 13641313.5
$ws.Range("M122").Value = -3112.52956
$ws.Range("N122").Value = -13646213.5
$ws.Range("H132").Value = 1251.25
$ws.Range("J132").Value = 1335
$ws.Range("L132").Value = 12015
$ws.Range("N132").Value = -17075
$ws.Range("H135").Value = 497
$ws.Range("I135").Value = 411.30435
$ws.Range("K135").Value = 3701.73915
$ws.Range("M135").Value = -1166.73915
$ws.Range("H139").Value = 3908.647
$ws.Range("I139").Value = 2736.6843
$ws.Range("J139").Value = 5393.1333
$ws.Range("K139").Value = 8210.052899999999
$ws.Range("L139").Value = 16179.3999
$ws.Range("M139").Value = -3070.052899999999
$ws.Range("N139").Value = -26459.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8842.125
$ws.Range("I68").Value = 13186
$ws.Range("J68").Value = 3257.1428
$ws.Range("K68").Value = 13186
$ws.Range("L68").Value = 3257.1428
$ws.Range("M68").Value = -12437
$ws.Range("N68").Value = -4755.1428
$ws.Range("H71").Value = 8842.125
$ws.Range("I71").Value = 13186
$ws.Range("J71").Value = 3257.1428
$ws.Range("K71").Value = 65930
$ws.Range("L71").Value = 16285.714
$ws.Range("M71").Value = -62186
$ws.Range("N71").Value = -23773.714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4459.0835
$ws.Range("I62").Value = 3900
$ws.Range("J62").Value = 4858.4287
$ws.Range("K62").Value = 3900
$ws.Range("L62").Value = 4858.4287
$ws.Range("M62").Value = -3276
$ws.Range("N62").Value = -6106.4287
$ws.Range("H65").Value = 4459.0835
$ws.Range("I65").Value = 3900
$ws.Range("J65").Value = 4858.4287
$ws.Range("K65").Value = 19500
$ws.Range("L65").Value = 24292.1435
$ws.Range("M65").Value = -16380
$ws.Range("N65").Value = -30532.1435
$ws.Range("H69").Value = 17317.75
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 22090.334
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 22090.334
$ws.Range("M69").Value = -2251
$ws.Range("N69").Value = -23588.334
$ws.Range("H72").Value = 17317.75
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 22090.334
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 66271.00199999999
$ws.Range("M72").Value = -5256
$ws.Range("N72").Value = -73759.00199999999
$ws.Range("H132").Value = 1054.5745
$ws.Range("I132").Value = 536.6896400000001
$ws.Range("J132").Value = 1888.9445
$ws.Range("K132").Value = 1610.06892
$ws.Range("L132").Value = 5666.833500000001
$ws.Range("M132").Value = 919.9310799999998
$ws.Range("N132").Value = -10726.8335
$ws.Range("H136").Value = 3270.2856
$ws.Range("I136").Value = 758.129
$ws.Range("K136").Value = 2274.387
$ws.Range("M136").Value = 275.6129999999998
